# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
#
# The source data rows (match records) for this league got re-synced and a
# few rows shifted position in the feed:
#   - Rows 352-357 (B:AC, i.e. everything except the running id in column A)
#     rotate up by one: 352<-353, 353<-354, 354<-355, 355<-356, 356<-357,
#     357<-352 (the old row-352 payload wraps around into row 357).
#   - Rows 425-426 (B:AC) simply swap with each other.
#
# Capture every source row's values BEFORE writing anything, then write the
# shifted/swapped data back out so earlier writes never clobber a value that
# is still needed later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 352-357, circular shift up by one row -------------------
$block1Rows = 352..357
$block1Data = @{}
foreach ($r in $block1Rows) {
    $block1Data[$r] = $ws.Range("B${r}:AC${r}").Value2
}

for ($i = 0; $i -lt $block1Rows.Length; $i++) {
    $destRow = $block1Rows[$i]
    $srcRow = $block1Rows[($i + 1) % $block1Rows.Length]
    $ws.Range("B${destRow}:AC${destRow}").Value2 = $block1Data[$srcRow]
}

# --- Block 2: rows 425-426, swap ---------------------------------------------
$block2Rows = 425..426
$block2Data = @{}
foreach ($r in $block2Rows) {
    $block2Data[$r] = $ws.Range("B${r}:AC${r}").Value2
}

for ($i = 0; $i -lt $block2Rows.Length; $i++) {
    $destRow = $block2Rows[$i]
    $srcRow = $block2Rows[($i + 1) % $block2Rows.Length]
    $ws.Range("B${destRow}:AC${destRow}").Value2 = $block2Data[$srcRow]
}
